$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 986, shifting the existing rows (986:1006) down to (989:1009).
$ws.Rows.Item(986).Resize(3).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# New row 986
$ws.Cells.Item(986, 1).Value = 10
$ws.Cells.Item(986, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(986, 3).Value = "La Araucanía"
$ws.Cells.Item(986, 4).Value = 45239
$ws.Cells.Item(986, 5).Value = 9
$ws.Cells.Item(986, 6).Value = 100112032
$ws.Cells.Item(986, 7).Value = "Zapallo italiano"
$ws.Cells.Item(986, 8).Value = "Huracán"
$ws.Cells.Item(986, 9).Value = "Primera"
$ws.Cells.Item(986, 10).Value = 100
$ws.Cells.Item(986, 11).Value = 14000
$ws.Cells.Item(986, 12).Value = 14000
$ws.Cells.Item(986, 13).Value = 14000
$ws.Cells.Item(986, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(986, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(986, 16).Value = 280
$ws.Cells.Item(986, 17).Value = 50
$ws.Cells.Item(986, 18).Value = "Hortaliza"

# New row 987
$ws.Cells.Item(987, 1).Value = 10
$ws.Cells.Item(987, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(987, 3).Value = "La Araucanía"
$ws.Cells.Item(987, 4).Value = 45239
$ws.Cells.Item(987, 5).Value = 9
$ws.Cells.Item(987, 6).Value = 100112032
$ws.Cells.Item(987, 7).Value = "Zapallo italiano"
$ws.Cells.Item(987, 8).Value = "Sin especificar"
$ws.Cells.Item(987, 9).Value = "Primera"
$ws.Cells.Item(987, 10).Value = 50
$ws.Cells.Item(987, 11).Value = 24000
$ws.Cells.Item(987, 12).Value = 24000
$ws.Cells.Item(987, 13).Value = 24000
$ws.Cells.Item(987, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(987, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(987, 16).Value = 480
$ws.Cells.Item(987, 17).Value = 50
$ws.Cells.Item(987, 18).Value = "Hortaliza"

# New row 988
$ws.Cells.Item(988, 1).Value = 10
$ws.Cells.Item(988, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(988, 3).Value = "La Araucanía"
$ws.Cells.Item(988, 4).Value = 45239
$ws.Cells.Item(988, 5).Value = 9
$ws.Cells.Item(988, 6).Value = 100112032
$ws.Cells.Item(988, 7).Value = "Zapallo italiano"
$ws.Cells.Item(988, 8).Value = "Sin especificar"
$ws.Cells.Item(988, 9).Value = "Primera"
$ws.Cells.Item(988, 10).Value = 50
$ws.Cells.Item(988, 11).Value = 20000
$ws.Cells.Item(988, 12).Value = 20000
$ws.Cells.Item(988, 13).Value = 20000
$ws.Cells.Item(988, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(988, 15).Value = "Región del Maule"
$ws.Cells.Item(988, 16).Value = 400
$ws.Cells.Item(988, 17).Value = 50
$ws.Cells.Item(988, 18).Value = "Hortaliza"
